$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.435.54'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.643.15'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.98'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('E6').Value = '  +4.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.21'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.257'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0609'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.02%  '
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').Value = '1.873.02'
$ws.Range('E12').Value = '  -1.52%  '
$ws.Range('D13').Value = '1.655.04'
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.03'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.30'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').Value = '27.382.06'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.34'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -8.37%  '
$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.49'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.31'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.03'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.53'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('E26').Value = '  +2.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.94'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.52'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -6.45%  '
$ws.Range('E30').Value = '  -4.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0487'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.27'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.09'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('D34').Value = '1.394.22'
$ws.Range('E34').Value = '  -5.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.56'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.881'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.559'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.23%  '
$ws.Range('E39').Value = '  -3.25%  '
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.47'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.17'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -8.00%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.785'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').Value = '1.784.09'
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.65'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.11'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0982'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.60%  '
